$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.607.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "'1.987.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.55%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'242.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  -3.31%  "
$ws.Range("D7").Value = "'57.39"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.38%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'60.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").Value = "'0.360"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "'0.0730"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("E12").Value = "  -4.76%  "
$ws.Range("D13").Value = "'0.922"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").Value = "'14.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.32%  "
$ws.Range("D15").Value = "'2.276.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.68%  "
$ws.Range("D16").Value = "'5.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "'1.987.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("D18").Value = "'17.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.79%  "
$ws.Range("D19").Value = "'35.523.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'70.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").Value = "'0.0₃0837"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("D22").Value = "'233.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("E23").Value = "  -3.21%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "'2.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.34%  "
$ws.Range("D27").Value = "'163.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'9.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").Value = "'19.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.42%  "
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -4.24%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "'0.0905"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.59%  "
$ws.Range("D35").Value = "'4.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.93%  "
$ws.Range("D36").Value = "'2.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("D42").Value = "'0.0210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").Value = "  -3.55%  "
$ws.Range("D44").Value = "'0.0889"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.42%  "
$ws.Range("D45").Value = "'91.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("D46").Value = "'1.380.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "'15.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'2.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.49%  "
$ws.Range("D51").Value = "'45.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.38%  "
